$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.55
$ws.Range("I2").Value = 2.9
$ws.Range("J2").Value = 3.25
$ws.Range("W2").Value = 8
$ws.Range("AL2").Value = 23
